# Update the "dSF" column (F) values to reflect repulled data / pushed
# data and the recalculated mean, per the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F5").Value = -5
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = -10
$ws.Range("F12").Value = -5
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = 2
$ws.Range("F17").Value = 4
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = -3
